# Append the new resale-number row (2024-01-14 11:47) as row 56,
# mirroring the existing data rows' layout: A-D text, E-T numeric.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 56

# Text columns, matching the existing rows (A-D are plain text, not dates
# or numbers). "2024-01-14" and "02" are at risk of being auto-coerced by
# Excel into a date serial / plain integer, so a leading apostrophe forces
# text interpretation for just those two. "11:47:20" and "Sunday" are not
# at risk, so they're set plainly (no text-coercion styling needed).
$ws.Cells.Item($row, 1).Value = "'2024-01-14"
$ws.Cells.Item($row, 2).Value = "11:47:20"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "'02"

# Numeric columns
$ws.Cells.Item($row, 5).Value  = 138964
$ws.Cells.Item($row, 6).Value  = 143003
$ws.Cells.Item($row, 7).Value  = 171054
$ws.Cells.Item($row, 8).Value  = 148326
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 119430
$ws.Cells.Item($row, 11).Value = 225027
$ws.Cells.Item($row, 12).Value = 253725
$ws.Cells.Item($row, 13).Value = 185248
$ws.Cells.Item($row, 14).Value = 110492
$ws.Cells.Item($row, 15).Value = 41053
$ws.Cells.Item($row, 16).Value = 30899
$ws.Cells.Item($row, 17).Value = 73172
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42511
$ws.Cells.Item($row, 20).Value = -1
